{"js": "// Rename the \"weight\" column to \"weight_kilogram\" (and \"maxWeight\" to\n// \"maxWeight_kilogram\") throughout the relational-schema table, and add an\n// extra trailing empty paragraph at the end of the document body.\n\nconst body = context.document.body;\n\n// 1) \"maxWeight\" -> \"maxWeight_kilogram\"  (Tray table header line)\nconst maxWeightResults = body.search(\"maxWeight\", { matchCase: true });\nmaxWeightResults.load(\"text\");\nawait context.sync();\n\nmaxWeightResults.items.forEach((item) => {\n  item.insertText(\"maxWeight_kilogram\", Word.InsertLocation.replace);\n});\nawait context.sync();\n\n// 2) \", weight\" -> \", weight_kilogram\"\n//    Occurs 3 times: Animal (\u2026, weight), AnimalPart (\u2026, weight, \u2026) and\n//    Tray (\u2026, weight, \u2026). Matching on \", weight\" (with the leading comma\n//    and space) keeps this from also matching the \"maxWeight\" substring.\nconst weightResults = body.search(\", weight\", { matchCase: true });\nweightResults.load(\"text\");\nawait context.sync();\n\nweightResults.items.forEach((item) => {\n  item.insertText(\", weight_kilogram\", Word.InsertLocation.replace);\n});\nawait context.sync();\n\n// 3) Add one extra empty paragraph right before the final (pre-existing)\n//    empty paragraph that precedes the section break, i.e. directly after\n//    the table.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"\", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Rename the \"weight\" column to \"weight_kilogram\" (and \"maxWeight\" to\n# \"maxWeight_kilogram\") throughout the relational-schema table, and add an\n# extra trailing empty paragraph at the end of the document body.\n\n$d = $word.ActiveDocument\n\n# 1) \"maxWeight\" -> \"maxWeight_kilogram\"  (Tray table header line)\n$find = $d.Content.Find\n$find.Execute(\"maxWeight\", $false, $false, $false, $false, $false, $true, 1, $false, \"maxWeight_kilogram\", 2)\n\n# 2) \", weight\" -> \", weight_kilogram\"\n#    Occurs 3 times: Animal (\u2026, weight), AnimalPart (\u2026, weight, \u2026) and\n#    Tray (\u2026, weight, \u2026). Matching on \", weight\" (with the leading comma\n#    and space) keeps this from also matching the \"maxWeight\" substring.\n#    wdReplaceAll (2) updates every occurrence in one call.\n$find2 = $d.Content.Find\n$find2.Execute(\", weight\", $false, $false, $false, $false, $false, $true, 1, $false, \", weight_kilogram\", 2)\n\n# 3) Add one extra empty paragraph right before the final (pre-existing)\n#    empty paragraph that precedes the section break, i.e. directly after\n#    the table.\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphBefore()\n"}
